$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 438; existing rows 438:531 shift down to 440:533.
$ws.Rows("438:439").Insert()

# New row 438: Navel Late / Primera, 2021-11-08
$ws.Range("A438").Value = 10
$ws.Range("B438").Value = "Vega Modelo de Temuco"
$ws.Range("C438").Value = "La Araucanía"
$ws.Range("D438").Value = 44508
$ws.Range("E438").Value = 9
$ws.Range("F438").Value = "Fruta"
$ws.Range("G438").Value = 100102
$ws.Range("H438").Value = "Cítricos"
$ws.Range("I438").Value = 100102005
$ws.Range("J438").Value = "Naranja"
$ws.Range("K438").Value = "Navel Late"
$ws.Range("L438").Value = "Primera"
$ws.Range("M438").Value = 155
$ws.Range("N438").Value = 10000
$ws.Range("O438").Value = 10000
$ws.Range("P438").Value = 10000
$ws.Range("Q438").Value = "`$/bandeja 15 kilos granel"
$ws.Range("R438").Value = "Región de O'Higgins"
$ws.Range("S438").Value = 667
$ws.Range("T438").Value = 15

# New row 439: Navel Late / Primera, 2021-11-08, bins
$ws.Range("A439").Value = 10
$ws.Range("B439").Value = "Vega Modelo de Temuco"
$ws.Range("C439").Value = "La Araucanía"
$ws.Range("D439").Value = 44508
$ws.Range("E439").Value = 9
$ws.Range("F439").Value = "Fruta"
$ws.Range("G439").Value = 100102
$ws.Range("H439").Value = "Cítricos"
$ws.Range("I439").Value = 100102005
$ws.Range("J439").Value = "Naranja"
$ws.Range("K439").Value = "Navel Late"
$ws.Range("L439").Value = "Primera"
$ws.Range("M439").Value = 20
$ws.Range("N439").Value = 170000
$ws.Range("O439").Value = 180000
$ws.Range("P439").Value = 177500
$ws.Range("Q439").Value = "`$/bins (400 kilos)"
$ws.Range("R439").Value = "Región de O'Higgins"
$ws.Range("S439").Value = 444
$ws.Range("T439").Value = 400
